$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("C1").Borders.Item(1).LineStyle = 1
Write-Output "done"
